$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# coinranking "Price" snapshot refresh (GitHub Actions cron).
# Price (D) and Volume(1h) (E) columns are always plain text in this sheet
# (prices use "."-grouped thousands, e.g. "26.700.57"; volumes are "  +0.20%  "
# style strings) - never real numbers/percentages. Some of the new Price
# values (e.g. "211.35") are lexically valid floats, so a plain .Value write
# would make Excel auto-coerce them into numbers. Guard those specific cells
# by forcing Text format right before the write, then drop the format change
# again with ClearFormats so the cell keeps its original (unstyled) look -
# only the stored type/content changes, exactly like the source diff.

$ws.Range("D2").Value = "26.700.57"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.599.60"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("E11").Value = "  +0.75%  "

$ws.Range("D12").Value = "1.824.19"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "1.614.39"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("D17").Value = "26.678.14"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "0.0₃0760"
$ws.Range("E18").Value = "  +4.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.22"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.66%  "

$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.92"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.02%  "

$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("E32").Value = "  +0.69%  "

$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("D34").Value = "1.290.27"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.621"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  +16.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.827"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.43"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.785"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.27"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("D45").Value = "1.733.02"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.94%  "

$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.38"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.68%  "

